# Regenerate the "K" column (column G) values for rows 2-28 of the save-data
# sheet. Per the commit message, the column previously held a "Strike#" style
# total that is being replaced with the true strikeout count (K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 4
    6  = 0
    7  = 3
    8  = 1
    9  = 3
    10 = 2
    11 = 1
    12 = 2
    13 = 6
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 3
    20 = 3
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 3
    26 = 3
    27 = 3
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
